# Auto-generated Excel COM-interop script
# Applies updated market-price / profit figures to the Maduin_Profits workbook
# (scheduled runner refresh of currentAveragePrice / LevePrice / LeveProfit columns)

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (17 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 10.285714
$ws.Range("I11").Value = 10.285714
$ws.Range("K11").Value = 10.285714
$ws.Range("M11").Value = 129.714286
$ws.Range("H33").Value = 435.91666
$ws.Range("I33").Value = 248.36363
$ws.Range("K33").Value = 248.36363
$ws.Range("M33").Value = -19.36363
$ws.Range("H64").Value = 12500
$ws.Range("H67").Value = 12500
$ws.Range("H132").Value = 2335.7144
$ws.Range("I132").Value = 2335.7144
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7007.1432
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -4477.1432

# --- Sheet: ARM (29 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3682.6875
$ws.Range("I32").Value = 3276.6
$ws.Range("K32").Value = 3276.6
$ws.Range("M32").Value = -2989.6
$ws.Range("H45").Value = 3335.8948
$ws.Range("I45").Value = 1230.3334
$ws.Range("K45").Value = 1230.3334
$ws.Range("M45").Value = -853.3334
$ws.Range("H88").Value = 2075
$ws.Range("I88").Value = 2037
$ws.Range("J88").Value = 2090.2
$ws.Range("K88").Value = 2037
$ws.Range("L88").Value = 2090.2
$ws.Range("M88").Value = -1631
$ws.Range("N88").Value = -2902.2
$ws.Range("H91").Value = 2075
$ws.Range("I91").Value = 2037
$ws.Range("J91").Value = 2090.2
$ws.Range("K91").Value = 2037
$ws.Range("L91").Value = 2090.2
$ws.Range("M91").Value = -633
$ws.Range("N91").Value = -4898.2
$ws.Range("H132").Value = 2154.818
$ws.Range("I132").Value = 2154.818
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6464.454000000001
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -3934.454000000001

# --- Sheet: BSM (15 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 2424.2
$ws.Range("I80").Value = 900
$ws.Range("K80").Value = 900
$ws.Range("M80").Value = 98
$ws.Range("H83").Value = 2424.2
$ws.Range("I83").Value = 900
$ws.Range("K83").Value = 4500
$ws.Range("M83").Value = 492
$ws.Range("H107").Value = 587.1667
$ws.Range("I107").Value = 467.81818
$ws.Range("J107").Value = 1900
$ws.Range("K107").Value = 467.81818
$ws.Range("L107").Value = 1900
$ws.Range("M107").Value = 1452.18182
$ws.Range("N107").Value = -5740

# --- Sheet: CRP (35 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 647.3333
$ws.Range("I58").Value = 476.8
$ws.Range("K58").Value = 476.8
$ws.Range("M58").Value = -273.8
$ws.Range("H94").Value = 3647.3635
$ws.Range("I94").Value = 4471.8
$ws.Range("J94").Value = 2960.3333
$ws.Range("K94").Value = 4471.8
$ws.Range("L94").Value = 2960.3333
$ws.Range("M94").Value = -4020.8
$ws.Range("N94").Value = -3862.3333
$ws.Range("H99").Value = 3743.8948
$ws.Range("I99").Value = 4039.625
$ws.Range("K99").Value = 4039.625
$ws.Range("M99").Value = -2541.625
$ws.Range("H106").Value = 61000.453
$ws.Range("J106").Value = 61000.453
$ws.Range("L106").Value = 61000.453
$ws.Range("N106").Value = -63524.453
$ws.Range("H122").Value = 2247.3333
$ws.Range("I122").Value = 2247.3333
$ws.Range("K122").Value = 6741.999899999999
$ws.Range("M122").Value = -4291.999899999999
$ws.Range("H126").Value = 3743.8948
$ws.Range("I126").Value = 4039.625
$ws.Range("K126").Value = 12118.875
$ws.Range("M126").Value = -9648.875
$ws.Range("H134").Value = 2181.25
$ws.Range("I134").Value = 2070.3333
$ws.Range("K134").Value = 6210.999899999999
$ws.Range("M134").Value = -3675.999899999999
$ws.Range("H136").Value = 647.3333
$ws.Range("I136").Value = 476.8
$ws.Range("K136").Value = 1430.4
$ws.Range("M136").Value = 1119.6

# --- Sheet: CUL (7 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 922.5
$ws.Range("I131").Value = 715
$ws.Range("J131").Value = 991.6667
$ws.Range("K131").Value = 2145
$ws.Range("L131").Value = 2975.0001
$ws.Range("M131").Value = 2895
$ws.Range("N131").Value = -13055.0001

# --- Sheet: GSM (22 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 25000
$ws.Range("J96").Value = 25000
$ws.Range("L96").Value = 25000
$ws.Range("N96").Value = -30492
$ws.Range("H113").Value = 10007.667
$ws.Range("I113").Value = 10008.5
$ws.Range("J113").Value = 10006
$ws.Range("K113").Value = 10008.5
$ws.Range("L113").Value = 10006
$ws.Range("M113").Value = -7838.5
$ws.Range("N113").Value = -14346
$ws.Range("H126").Value = 4841.933
$ws.Range("I126").Value = 4077.9
$ws.Range("K126").Value = 12233.7
$ws.Range("M126").Value = -9763.700000000001
$ws.Range("H132").Value = 2957.6
$ws.Range("I132").Value = 2957.6
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8872.799999999999
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -6342.799999999999

# --- Sheet: LTW (59 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 5010
$ws.Range("I4").Value = 5010
$ws.Range("K4").Value = 5010
$ws.Range("M4").Value = -4897
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H22").Value = 6944.222
$ws.Range("J22").Value = 7749.75
$ws.Range("L22").Value = 7749.75
$ws.Range("N22").Value = -8339.75
$ws.Range("H27").Value = 6944.222
$ws.Range("J27").Value = 7749.75
$ws.Range("L27").Value = 7749.75
$ws.Range("N27").Value = -7963.75
$ws.Range("H28").Value = 5010
$ws.Range("I28").Value = 5010
$ws.Range("K28").Value = 5010
$ws.Range("M28").Value = -4778
$ws.Range("H37").Value = 5010
$ws.Range("I37").Value = 5010
$ws.Range("K37").Value = 5010
$ws.Range("M37").Value = -4903
$ws.Range("H38").Value = 40000
$ws.Range("J38").Value = 40000
$ws.Range("L38").Value = 40000
$ws.Range("N38").Value = -40820
$ws.Range("H40").Value = 1769.3334
$ws.Range("I40").Value = 1769.3334
$ws.Range("K40").Value = 1769.3334
$ws.Range("M40").Value = -1633.3334
$ws.Range("H47").Value = 20000
$ws.Range("J47").Value = 20000
$ws.Range("L47").Value = 20000
$ws.Range("N47").Value = -20980
$ws.Range("H52").Value = 20000
$ws.Range("J52").Value = 20000
$ws.Range("L52").Value = 20000
$ws.Range("N52").Value = -20466
$ws.Range("H82").Value = 430.6
$ws.Range("I82").Value = 425.75
$ws.Range("K82").Value = 425.75
$ws.Range("M82").Value = -64.75
$ws.Range("H85").Value = 430.6
$ws.Range("I85").Value = 425.75
$ws.Range("K85").Value = 425.75
$ws.Range("M85").Value = 822.25
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H132").Value = 3224.9
$ws.Range("I132").Value = 2530.5
$ws.Range("J132").Value = 6002.5
$ws.Range("K132").Value = 7591.5
$ws.Range("L132").Value = 18007.5
$ws.Range("M132").Value = -5061.5
$ws.Range("N132").Value = -23067.5

# --- Sheet: WVR (24 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 8000
$ws.Range("I7").Value = 8000
$ws.Range("K7").Value = 8000
$ws.Range("M7").Value = -7887
$ws.Range("H39").Value = 10000
$ws.Range("I39").Value = 10000
$ws.Range("K39").Value = 10000
$ws.Range("M39").Value = -9587
$ws.Range("I81").Value = 5511.3335
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 11022.667
$ws.Range("L81").Value = 4000
$ws.Range("M81").Value = -9961.666999999999
$ws.Range("N81").Value = -6122
$ws.Range("I84").Value = 5511.3335
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 55113.335
$ws.Range("L84").Value = 20000
$ws.Range("M84").Value = -49809.335
$ws.Range("N84").Value = -30608
$ws.Range("H126").Value = 1860.7858
$ws.Range("I126").Value = 2521.8572
$ws.Range("K126").Value = 7565.571599999999
$ws.Range("M126").Value = -5095.571599999999
